$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------
# Part 1: the underlying match-odds data for these rows was re-fetched
# and now sorts differently; the row content (columns B..AC) for rows
# 130-136 and 142-145 gets reshuffled among themselves. Column A (the
# running index) stays fixed per row. Capture the current B:AC values
# for every affected row first, then write them back out in their new
# positions, so overlapping/cyclic reassignments don't clobber data
# that still needs to be read.
# -------------------------------------------------------------------

$snapshot = @{}
$affectedRows = @(130,131,132,133,134,135,136,142,143,144,145)
foreach ($r in $affectedRows) {
    $snapshot[$r] = $ws.Range("B$r`:AC$r").Value2
}

# target row -> source row whose old content it should receive
$mapping = @{
    130 = 133
    131 = 132
    132 = 131
    133 = 130
    134 = 135
    135 = 136
    136 = 134
    142 = 145
    143 = 142
    144 = 143
    145 = 144
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $ws.Range("B$targetRow`:AC$targetRow").Value = $snapshot[$sourceRow]
}

# -------------------------------------------------------------------
# Part 2: append 4 newly scraped matches as rows 212-215.
# -------------------------------------------------------------------

# Copy formatting (styles) from the last existing row so the new rows
# pick up the same column A / column E number formats.
$ws.Range("A211:AC211").Copy() | Out-Null
$ws.Range("A212:AC215").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row=212; A=210; B=7773508; C="Ecuador LigaPro Serie A"; D="Ecuador LigaPro Serie A"; E=45400.875;          F="Barcelona Guayaquil"; G="El Nacional";   H=1; I=0; J="H"; K=1.45;  L=4.5;  M=6.5;  N=1.45;  O=4.5;   P=7;    Q=-1.25; R=2;     S=1.8;   T=2.75; U=1.925; V=1.875; W=0.45;  X=-1;   Y=-1; Z=-0.5;  AA=0.4;    AB=-1;   AC=0.875 },
    @{ Row=213; A=211; B=7773509; C="Ecuador LigaPro Serie A"; D="Ecuador LigaPro Serie A"; E=45401.66666666666;  F="Delfin SC";           G="Macara";        H=1; I=1; J="D"; K=2;     L=3.25; M=3.75; N=2.25;  O=3.1;   P=3.25; Q=-0.25; R=1.925; S=1.875; T=2;    U=1.8;   V=2;     W=-1;    X=2.1;  Y=-1; Z=-0.5;  AA=0.4375; AB=0;    AC=0     },
    @{ Row=214; A=212; B=7773506; C="Ecuador LigaPro Serie A"; D="Ecuador LigaPro Serie A"; E=45401.77083333334;  F="SD Aucas";            G="Emelec";        H=3; I=1; J="H"; K=2;     L=3.25; M=3.6;  N=1.8;   O=3.4;   P=4.2;  Q=-0.5;  R=1.825; S=1.975; T=2.25; U=1.85;  V=1.95;  W=0.8;   X=-1;   Y=-1; Z=0.825; AA=-1;     AB=0.85; AC=-1    },
    @{ Row=215; A=213; B=7773507; C="Ecuador LigaPro Serie A"; D="Ecuador LigaPro Serie A"; E=45401.77083333334;  F="Orense";              G="LDU Quito";     H=1; I=0; J="H"; K=4;     L=3.25; M=1.909; N=3.4;   O=3;     P=2.2;  Q=0.25;  R=1.9;   S=1.9;   T=2.5;  U=1.975; V=1.825; W=2.4;   X=-1;   Y=-1; Z=0.9;   AA=-1;     AB=-1;   AC=0.825 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Range("A$r").Value = $nr.A
    $ws.Range("B$r").Value = $nr.B
    $ws.Range("C$r").Value = $nr.C
    $ws.Range("D$r").Value = $nr.D
    $ws.Range("E$r").Value = $nr.E
    $ws.Range("F$r").Value = $nr.F
    $ws.Range("G$r").Value = $nr.G
    $ws.Range("H$r").Value = $nr.H
    $ws.Range("I$r").Value = $nr.I
    $ws.Range("J$r").Value = $nr.J
    $ws.Range("K$r").Value = $nr.K
    $ws.Range("L$r").Value = $nr.L
    $ws.Range("M$r").Value = $nr.M
    $ws.Range("N$r").Value = $nr.N
    $ws.Range("O$r").Value = $nr.O
    $ws.Range("P$r").Value = $nr.P
    $ws.Range("Q$r").Value = $nr.Q
    $ws.Range("R$r").Value = $nr.R
    $ws.Range("S$r").Value = $nr.S
    $ws.Range("T$r").Value = $nr.T
    $ws.Range("U$r").Value = $nr.U
    $ws.Range("V$r").Value = $nr.V
    $ws.Range("W$r").Value = $nr.W
    $ws.Range("X$r").Value = $nr.X
    $ws.Range("Y$r").Value = $nr.Y
    $ws.Range("Z$r").Value = $nr.Z
    $ws.Range("AA$r").Value = $nr.AA
    $ws.Range("AB$r").Value = $nr.AB
    $ws.Range("AC$r").Value = $nr.AC
}
